$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 7 (old) : "Dennis Audu, Orang Tang Enow, Mohammadsadegh Firouzi" ---
$p7xml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:t xml:space="preserve">Dennis </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Audu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Orang Tang E</w:t></w:r><w:r><w:t xml:space="preserve">now, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mohammadsadegh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Firouzi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(7).Range.InsertXML($p7xml)

# --- Paragraph 6 (old) : "Xin Zhao, Neel Ajay Mahimkar" ---
$p6xml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:t xml:space="preserve">Xin Zhao, </w:t></w:r><w:r><w:t xml:space="preserve">Neel Ajay </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mahimkar</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(6).Range.InsertXML($p6xml)

# --- Paragraph 5 (old) : "Huu Minh Phong Nguyen, Abhi Nileshkumar Patel" (fr-FR) ---
# Splits into two paragraphs: "Huu Minh Phong Nguyen" and a new "Abhi Nileshkumar Patel"
$p5xml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:contextualSpacing w:val="0"/><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Huu Minh </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Phong</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> Nguyen</w:t></w:r></w:p><w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:contextualSpacing w:val="0"/><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Abhi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Nileshkumar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> Patel</w:t></w:r></w:p>'
$d.Paragraphs(5).Range.InsertXML($p5xml)

# --- Paragraph 4 (old) : "Yiyuan Dong" ---
$p4xml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:contextualSpacing w:val="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Yiyuan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Dong</w:t></w:r></w:p>'
$d.Paragraphs(4).Range.InsertXML($p4xml)

# --- Paragraph 3 (old) : "Bhavjot Pal, Samay Sehgal, Kannav Sethi" ---
$p3xml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:contextualSpacing w:val="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Bhavjot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Pal, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Samay</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Sehgal, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Kannav</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Sethi</w:t></w:r></w:p>'
$d.Paragraphs(3).Range.InsertXML($p3xml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
